# Apply the "Added columns for path parameters" change:
#  - Tests sheet: insert a new "param:scope" column before "param:name",
#    populate it for the existing row, and append a new test row that
#    exercises the missing-required-parameter case.
#  - Documentation sheet: insert a "param:scope" row in the parameter
#    description list, and a "Required parameters" note near the bottom.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Tests"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tests")

# Insert a new column I ("param:scope"), shifting the existing
# param:* columns (old I..N) one slot to the right (new J..O).
$ws.Columns("I:I").Insert()

# Match the width of the other param:* columns (20 "characters"),
# using an existing column as the source so the stored width matches
# exactly instead of drifting due to unit-conversion rounding.
$ws.Range("I1").ColumnWidth = $ws.Range("J1").ColumnWidth

# Header for the new column.
$ws.Range("I1").Value = "param:scope"

# Existing data row gets a value for the new column.
$ws.Range("I2").Value = "work"

# New row 3: a validation test covering a missing required parameter.
# Values that look like booleans/numbers ("true", "1", "20", and even the
# blank param value) are prefixed with a leading apostrophe so Excel keeps
# them as literal text instead of being auto-coerced into boolean/numeric
# cells; ClearFormats() afterwards drops the "quote prefix" marker Excel
# otherwise records for those cells, so they come out as plain text cells
# like their counterparts in row 2.
$ws.Range("A3").Value = "get-related-list - Missing Required Param"
$ws.Range("B3").Value = "Test GET /api/related-list/:scope with missing required parameters"
$ws.Range("C3").Value = "'true"
$ws.Range("D3").Value = 400
$ws.Range("E3").Value = 10000
$ws.Range("F3").Value = 2000
$ws.Range("G3").Value = 500
$ws.Range("H3").Value = "get-related-list,validation"
$ws.Range("I3").Value = "'"
$ws.Range("J3").Value = "classification"
$ws.Range("K3").Value = "https://lux.collections.yale.edu/data/test/example"
$ws.Range("L3").Value = "'1"
$ws.Range("M3").Value = "'20"
$ws.Range("N3").Value = "'true"
$ws.Range("O3").Value = "'1"
$ws.Range("C3").ClearFormats()
$ws.Range("I3").ClearFormats()
$ws.Range("L3").ClearFormats()
$ws.Range("M3").ClearFormats()
$ws.Range("N3").ClearFormats()
$ws.Range("O3").ClearFormats()

# ---------------------------------------------------------------------
# Sheet 2: "Documentation"
# ---------------------------------------------------------------------
$doc = $wb.Worksheets.Item("Documentation")

# Insert the new "param:scope" description row right before "param:name".
$doc.Rows("18:18").Insert()
$doc.Range("A18").Value = "param:scope"
$doc.Range("B18").Value = "Search scope (work, person, place, concept, event, etc.) (string) (REQUIRED - highlighted in yellow)"

# Insert a row for the "Required parameters" note, right before the
# existing "Optional parameters" note (now at row 30 after the shift above).
$doc.Rows("30:30").Insert()
$doc.Range("A30").Value = "• Required parameters: scope"
